$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# Fill in the Day 9 summary values (rows 45-47)
$ws.Range("C45").Value = 7050
$ws.Range("C46").Value = 2430
$ws.Range("C47").Value = 2430

# Update the active selection on the sheet
$ws.Range("C47").Select()
